$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same across all data rows 2-26), columns B..Q
$values = @{
    "B" = 0.3082271070678201
    "C" = -2.241727458159763
    "D" = -11.23292888890734
    "E" = 0.1885856136292009
    "F" = -0.1356755099388904
    "G" = 0.4106662722775318
    "H" = 1.924429454498304
    "I" = 0.1582508707211246
    "J" = 0.3500857381124779
    "K" = 0.2541683044168013
    "L" = 0.2690475247287556
    "M" = 0.6408324837877148
    "N" = -1.07531867879654
    "O" = 0.6681140493286761
    "P" = 33.77994876784079
    "Q" = 53.281961965732
}

foreach ($col in $values.Keys) {
    $rng = $ws.Range($col + "2:" + $col + "26")
    $rng.Value = $values[$col]
}
